$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original styles for the Price/Volume columns, then force
# a Text number format while we write the new values so that numeric-
# looking strings (e.g. "27.102.24", "0.00001084") are not silently
# reinterpreted by Excel as numbers. Restore the original style
# afterwards so the workbook formatting is left untouched.
$origD = $ws.Range("D2:D51").Style
$origE = $ws.Range("E2:E51").Style
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.102.24'
$ws.Range('E2').Value = '  -1.24%  '

$ws.Range('D3').Value = '1.781.72'
$ws.Range('E3').Value = '  -2.17%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '337.03'
$ws.Range('E5').Value = '  -2.52%  '

$ws.Range('E6').Value = '  +0.05%  '

$ws.Range('D7').Value = '0.3816'
$ws.Range('E7').Value = '  -0.40%  '

$ws.Range('D8').Value = '0.3408'
$ws.Range('E8').Value = '  -3.46%  '

$ws.Range('D9').Value = '48.06'
$ws.Range('E9').Value = '  -2.93%  '

$ws.Range('D10').Value = '1.185'
$ws.Range('E10').Value = '  -4.39%  '

$ws.Range('D11').Value = '0.07433'
$ws.Range('E11').Value = '  -4.77%  '

$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  -0.10%  '

$ws.Range('D13').Value = '21.61'
$ws.Range('E13').Value = '  -2.84%  '

$ws.Range('D14').Value = '6.444'
$ws.Range('E14').Value = '  -3.28%  '

$ws.Range('D15').Value = '1.778.19'
$ws.Range('E15').Value = '  -2.48%  '

$ws.Range('D16').Value = '7.059'
$ws.Range('E16').Value = '  -3.00%  '

$ws.Range('D17').Value = '0.00001084'
$ws.Range('E17').Value = '  -4.06%  '

$ws.Range('D18').Value = '0.06634'
$ws.Range('E18').Value = '  -1.81%  '

$ws.Range('D19').Value = '83.32'
$ws.Range('E19').Value = '  -3.94%  '

$ws.Range('E20').Value = '  +0.11%  '

$ws.Range('D21').Value = '6.535'
$ws.Range('E21').Value = '  -0.72%  '

$ws.Range('D22').Value = '17.30'
$ws.Range('E22').Value = '  -2.68%  '

$ws.Range('D23').Value = '27.104.07'
$ws.Range('E23').Value = '  -1.33%  '

$ws.Range('D24').Value = '12.21'
$ws.Range('E24').Value = '  -8.36%  '

$ws.Range('D25').Value = '2.370'
$ws.Range('E25').Value = '  -3.69%  '

$ws.Range('D26').Value = '2.501'
$ws.Range('E26').Value = '  -7.47%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '21.03'
$ws.Range('E27').Value = '  -5.76%  '

$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = '1.453'
$ws.Range('E28').Value = '  -3.91%  '

$ws.Range('D29').Value = '155.10'
$ws.Range('E29').Value = '  +0.67%  '

$ws.Range('D30').Value = '1.982.35'
$ws.Range('E30').Value = '  -2.18%  '

$ws.Range('D31').Value = '133.75'
$ws.Range('E31').Value = '  -2.49%  '

$ws.Range('D32').Value = '3.982'
$ws.Range('E32').Value = '  -2.17%  '

$ws.Range('D33').Value = '6.019'
$ws.Range('E33').Value = '  -6.03%  '

$ws.Range('D34').Value = '0.08673'
$ws.Range('E34').Value = '  -1.69%  '

$ws.Range('D35').Value = '13.08'
$ws.Range('E35').Value = '  -7.65%  '

$ws.Range('D36').Value = '1.625'
$ws.Range('E36').Value = '  -4.71%  '

$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '5.383'
$ws.Range('E37').Value = '  -5.34%  '

$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D38').Value = '0.6822'
$ws.Range('E38').Value = '  -3.85%  '

$ws.Range('D39').Value = '0.06274'
$ws.Range('E39').Value = '  -4.48%  '

$ws.Range('D40').Value = '0.2171'
$ws.Range('E40').Value = '  -5.19%  '

$ws.Range('D41').Value = '0.02315'
$ws.Range('E41').Value = '  -4.68%  '

$ws.Range('D42').Value = '8.533'
$ws.Range('E42').Value = '  -5.79%  '

$ws.Range('D43').Value = '1.235'
$ws.Range('E43').Value = '  -5.24%  '

$ws.Range('D44').Value = '14.21'
$ws.Range('E44').Value = '  -4.98%  '

$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.07%  '

$ws.Range('D46').Value = '0.6421'
$ws.Range('E46').Value = '  -3.07%  '

$ws.Range('D47').Value = '3.857'
$ws.Range('E47').Value = '  -4.66%  '

$ws.Range('D48').Value = '2.120'
$ws.Range('E48').Value = '  -3.53%  '

$ws.Range('D49').Value = '131.23'
$ws.Range('E49').Value = '  -1.51%  '

$ws.Range('D50').Value = '0.07082'
$ws.Range('E50').Value = '  -3.88%  '

$ws.Range('D51').Value = '78.64'
$ws.Range('E51').Value = '  -3.07%  '

# Restore original number formats / styles.
$ws.Range("D2:D51").Style = $origD
$ws.Range("E2:E51").Style = $origE
